$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update validation output values (I2:I4) - AVERAGEIF formula in I6 will recalc automatically
$ws.Range("I2").Value = 89.501999999999995
$ws.Range("I3").Value = 89.501999999999995
$ws.Range("I4").Value = 89.501999999999995

# Update view: scroll to show column E and move selection to L8
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L8").Select()
